$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the outdated "No. of Sites ..." breakdown columns (X:AB) and the
# DIFFERENCE column (AD) for row 2, keeping PREVIOUS ACCOMPLISHMENT (AC2).
$ws.Range("X2:AB2").ClearContents()
$ws.Range("AD2").ClearContents()
